# Daily attendance processing - 2025-12-28 17:00:34
#
# Normalize the "Recorded By" column (G) on the "Session Analysis Results"
# sheet: whenever the automated "System" recorder appears together with
# other recorders in the comma-separated list, move "System" to the front
# of the list while preserving the relative order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# Column G is the 7th column ("Recorded By")
$col = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $value = $cell.Value2

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ",\s*"
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    # Use case-sensitive .Equals() comparisons: the data can contain both
    # "system" and "System" as distinct entries, and only the exact
    # "System" entry (capital S) should be relocated to the front.
    $hasExactSystem = $false
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) {
            $hasExactSystem = $true
        }
    }

    if ($hasExactSystem) {
        $rest = @()
        $removedOne = $false
        foreach ($p in $trimmed) {
            if ((-not $removedOne) -and $p.Equals("System")) {
                $removedOne = $true
                continue
            }
            $rest += $p
        }
        $newParts = @("System") + $rest
        $newValue = $newParts -join ", "

        if (-not $newValue.Equals($value)) {
            $cell.Value2 = $newValue
        }
    }
}
